$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 371 (existing rows 371+ shift down by 2)
$ws.Rows.Item(371).Insert()
$ws.Rows.Item(371).Insert()

# Populate the new row 371 (Primera)
$ws.Cells.Item(371, 1).Value = 3
$ws.Cells.Item(371, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(371, 3).Value = "Coquimbo"
$ws.Cells.Item(371, 4).Value = 44522
$ws.Cells.Item(371, 5).Value = 5
$ws.Cells.Item(371, 6).Value = 100112008
$ws.Cells.Item(371, 7).Value = "Coliflor"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 2800
$ws.Cells.Item(371, 11).Value = 550
$ws.Cells.Item(371, 12).Value = 600
$ws.Cells.Item(371, 13).Value = 573
$ws.Cells.Item(371, 14).Value = "$/unidad"
$ws.Cells.Item(371, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(371, 16).Value = 573
$ws.Cells.Item(371, 17).Value = 1
$ws.Cells.Item(371, 18).Value = "Hortaliza"

# Populate the new row 372 (Segunda)
$ws.Cells.Item(372, 1).Value = 3
$ws.Cells.Item(372, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44522
$ws.Cells.Item(372, 5).Value = 5
$ws.Cells.Item(372, 6).Value = 100112008
$ws.Cells.Item(372, 7).Value = "Coliflor"
$ws.Cells.Item(372, 8).Value = "Sin especificar"
$ws.Cells.Item(372, 9).Value = "Segunda"
$ws.Cells.Item(372, 10).Value = 1200
$ws.Cells.Item(372, 11).Value = 450
$ws.Cells.Item(372, 12).Value = 450
$ws.Cells.Item(372, 13).Value = 450
$ws.Cells.Item(372, 14).Value = "$/unidad"
$ws.Cells.Item(372, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(372, 16).Value = 450
$ws.Cells.Item(372, 17).Value = 1
$ws.Cells.Item(372, 18).Value = "Hortaliza"
